# Fruta / hortaliza, semanal
# Insert 3 new weekly price records (for date 2022-06-02 / serial 44714) at the
# top of the Brocoli data block (rows 405-407), pushing all existing records
# in that block down by 3 rows (old 405..429 -> new 408..432). Sheet
# dimension grows from A1:R429 to A1:R432.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 405:407 - Excel shifts rows 405-429 down to 408-432
# and the new blank rows inherit the neighbouring row's formatting (keeps the
# date-format style on column D).
$ws.Rows("405:407").Insert()

# Row 405 - Provincia del Elquí
$ws.Cells.Item(405, 1).Value  = 10
$ws.Cells.Item(405, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(405, 3).Value  = "La Araucanía"
$ws.Cells.Item(405, 4).Value  = 44714
$ws.Cells.Item(405, 5).Value  = 9
$ws.Cells.Item(405, 6).Value  = 100112023
$ws.Cells.Item(405, 7).Value  = "Brócoli"
$ws.Cells.Item(405, 8).Value  = "Sin especificar"
$ws.Cells.Item(405, 9).Value  = "Primera"
$ws.Cells.Item(405, 10).Value = 450
$ws.Cells.Item(405, 11).Value = 1300
$ws.Cells.Item(405, 12).Value = 1300
$ws.Cells.Item(405, 13).Value = 1300
$ws.Cells.Item(405, 14).Value = "`$/unidad"
$ws.Cells.Item(405, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(405, 16).Value = 1300
$ws.Cells.Item(405, 17).Value = 1
$ws.Cells.Item(405, 18).Value = "Hortaliza"

# Row 406 - Región Metropolitana
$ws.Cells.Item(406, 1).Value  = 10
$ws.Cells.Item(406, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(406, 3).Value  = "La Araucanía"
$ws.Cells.Item(406, 4).Value  = 44714
$ws.Cells.Item(406, 5).Value  = 9
$ws.Cells.Item(406, 6).Value  = 100112023
$ws.Cells.Item(406, 7).Value  = "Brócoli"
$ws.Cells.Item(406, 8).Value  = "Sin especificar"
$ws.Cells.Item(406, 9).Value  = "Primera"
$ws.Cells.Item(406, 10).Value = 1150
$ws.Cells.Item(406, 11).Value = 1000
$ws.Cells.Item(406, 12).Value = 1200
$ws.Cells.Item(406, 13).Value = 1096
$ws.Cells.Item(406, 14).Value = "`$/unidad"
$ws.Cells.Item(406, 15).Value = "Región Metropolitana"
$ws.Cells.Item(406, 16).Value = 1096
$ws.Cells.Item(406, 17).Value = 1
$ws.Cells.Item(406, 18).Value = "Hortaliza"

# Row 407 - Región del Maule
$ws.Cells.Item(407, 1).Value  = 10
$ws.Cells.Item(407, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(407, 3).Value  = "La Araucanía"
$ws.Cells.Item(407, 4).Value  = 44714
$ws.Cells.Item(407, 5).Value  = 9
$ws.Cells.Item(407, 6).Value  = 100112023
$ws.Cells.Item(407, 7).Value  = "Brócoli"
$ws.Cells.Item(407, 8).Value  = "Sin especificar"
$ws.Cells.Item(407, 9).Value  = "Primera"
$ws.Cells.Item(407, 10).Value = 850
$ws.Cells.Item(407, 11).Value = 1000
$ws.Cells.Item(407, 12).Value = 1000
$ws.Cells.Item(407, 13).Value = 1000
$ws.Cells.Item(407, 14).Value = "`$/unidad"
$ws.Cells.Item(407, 15).Value = "Región del Maule"
$ws.Cells.Item(407, 16).Value = 1000
$ws.Cells.Item(407, 17).Value = 1
$ws.Cells.Item(407, 18).Value = "Hortaliza"
